$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shrink the saved window height in the workbook view
$excel.ActiveWindow.Height = 8760

# Widen column A slightly
$ws.Columns.Item(1).ColumnWidth = 12

# New header row: First Name / Last Name / email
$ws.Range("B1").Value = "First Name"
$ws.Range("C1").Value = "Last Name"
$ws.Range("D1").Value = "email"

# New data row 2 values
$ws.Range("B2").Value = "Samyuktha"
$ws.Range("C2").Value = "CS"
$ws.Range("D2").Value = "abc@gmail.com"

# Turn the email address into a real mailto hyperlink, and apply the
# built-in "Hyperlink" style to match Excel's default behaviour
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:abc@gmail.com", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "abc@gmail.com")

# Update the existing number in A3
$ws.Range("A3").Value = 76755645454

Write-Host "edit applied"
